$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the bold header style from the existing E1 header cell onto the
# new F1:G1 header cells.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# Populate Password column first (G), then Username column (F), matching
# the order the data was entered.
$ws.Range("G1").Value = "Password"
$ws.Range("G2").Value = "nivi"
$ws.Range("G3").Value = "rahul"
$ws.Range("G4").Value = "tanu"

$ws.Range("F1").Value = "Username"
$ws.Range("F2").Value = "Nivethetha"
$ws.Range("F3").Value = "Rahul"
$ws.Range("F4").Value = "Tanu"

# Update Balance values
$ws.Range("C2").Value = 15704
$ws.Range("C3").Value = 20000

# Leave the selection where the author left it when saving the file
$ws.Range("H16").Select() | Out-Null
